$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "63.277.08"
$ws.Range("E2").Value = "  +1.30%  "

Set-TextValue "D3" "2.453.91"
$ws.Range("E3").Value = "  +1.14%  "

Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue "D5" "573.20"
$ws.Range("E5").Value = "  +1.33%  "

Set-TextValue "D6" "146.44"
$ws.Range("E6").Value = "  +0.99%  "

$ws.Range("E7").Value = "  +0.08%  "

Set-TextValue "D8" "0.538"
$ws.Range("E8").Value = "  +1.21%  "

Set-TextValue "D9" "2.451.24"
$ws.Range("E9").Value = "  +1.09%  "

$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("E11").Value = "  +1.32%  "

Set-TextValue "D12" "5.27"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("E13").Value = "  +0.40%  "

Set-TextValue "D14" "27.06"
$ws.Range("E14").Value = "  +1.16%  "

Set-TextValue "D15" "0.0000180"
$ws.Range("E15").Value = "  +0.68%  "

Set-TextValue "D16" "2.901.45"
$ws.Range("E16").Value = "  +3.86%  "

Set-TextValue "D17" "63.288.19"
$ws.Range("E17").Value = "  +1.64%  "

Set-TextValue "D18" "2.446.08"
$ws.Range("E18").Value = "  +1.32%  "

Set-TextValue "D19" "11.34"
$ws.Range("E19").Value = "  +1.27%  "

Set-TextValue "D20" "7.31"
$ws.Range("E20").Value = "  +5.21%  "

Set-TextValue "D21" "329.28"
$ws.Range("E21").Value = "  +1.97%  "

Set-TextValue "D22" "4.21"
$ws.Range("E22").Value = "  +1.32%  "

Set-TextValue "D23" "2.08"
$ws.Range("E23").Value = "  +14.90%  "

Set-TextValue "D24" "0.999"
$ws.Range("E24").Value = "  -0.06%  "

Set-TextValue "D25" "65.56"
$ws.Range("E25").Value = "  -2.11%  "

Set-TextValue "D26" "615.77"
$ws.Range("E26").Value = "  +3.95%  "

Set-TextValue "D27" "8.86"
$ws.Range("E27").Value = "  +3.70%  "

$ws.Range("E28").Value = "  +3.21%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D29" "2.575.80"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D30" "1.51"
$ws.Range("E30").Value = "  +4.81%  "

Set-TextValue "D31" "0.995"
$ws.Range("E31").Value = "  -0.43%  "

Set-TextValue "D32" "8.23"
$ws.Range("E32").Value = "  -2.05%  "

$ws.Range("E33").Value = "  +1.66%  "

$ws.Range("E34").Value = "  -2.00%  "

Set-TextValue "D35" "5.20"
$ws.Range("E35").Value = "  +7.58%  "

$ws.Range("E36").Value = "  +1.95%  "

Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D39" "18.87"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D40" "5.42"
$ws.Range("E40").Value = "  +1.51%  "

Set-TextValue "D41" "146.76"
$ws.Range("E41").Value = "  -0.33%  "

Set-TextValue "D42" "1.79"
$ws.Range("E42").Value = "  -1.19%  "

Set-TextValue "D43" "2.61"
$ws.Range("E43").Value = "  +6.74%  "

$ws.Range("E44").Value = "  -0.03%  "

Set-TextValue "D45" "41.79"
$ws.Range("E45").Value = "  +0.51%  "

Set-TextValue "D46" "148.67"
$ws.Range("E46").Value = "  +0.56%  "

Set-TextValue "D47" "3.77"
$ws.Range("E47").Value = "  +3.11%  "

Set-TextValue "D48" "21.17"
$ws.Range("E48").Value = "  +3.52%  "

Set-TextValue "D49" "0.0534"
$ws.Range("E49").Value = "  -0.01%  "

Set-TextValue "D50" "0.602"
$ws.Range("E50").Value = "  +0.26%  "

Set-TextValue "D51" "0.0233"
$ws.Range("E51").Value = "  +1.12%  "
